# Scheduled-runner style refresh of market-price / profit columns (H:N)
# across the profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below were re-pulled from the market data source; only the
# price/profit columns for specific leve rows are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1242.2941
$ws.Range("I28").Value = 1195.8889
$ws.Range("J28").Value = 1294.5
$ws.Range("K28").Value = 1195.8889
$ws.Range("L28").Value = 1294.5
$ws.Range("M28").Value = -710.8888999999999
$ws.Range("N28").Value = -2264.5

$ws.Range("H40").Value = 1848
$ws.Range("I40").Value = 2660
$ws.Range("J40").Value = 1645
$ws.Range("K40").Value = 2660
$ws.Range("L40").Value = 1645
$ws.Range("M40").Value = -2485
$ws.Range("N40").Value = -1995

$ws.Range("H98").Value = 731.17645
$ws.Range("I98").Value = 761.6
$ws.Range("K98").Value = 761.6
$ws.Range("M98").Value = 736.4

$ws.Range("H106").Value = 4850.357
$ws.Range("I106").Value = 5050.4165
$ws.Range("K106").Value = 5050.4165
$ws.Range("M106").Value = -4419.4165

$ws.Range("H107").Value = 400
$ws.Range("I107").Value = 400
$ws.Range("K107").Value = 400
$ws.Range("M107").Value = 1520

$ws.Range("H112").Value = 29413400
$ws.Range("I112").Value = 250000240
$ws.Range("J112").Value = 1820.3
$ws.Range("K112").Value = 750000720
$ws.Range("L112").Value = 5460.9
$ws.Range("M112").Value = -749999612
$ws.Range("N112").Value = -7676.9

$ws.Range("H116").Value = 1916.375
$ws.Range("I116").Value = 1768.75
$ws.Range("K116").Value = 1768.75
$ws.Range("M116").Value = 1673.25

$ws.Range("H122").Value = 731.17645
$ws.Range("I122").Value = 761.6
$ws.Range("K122").Value = 2284.8
$ws.Range("M122").Value = 165.1999999999998

$ws.Range("H123").Value = 38000
$ws.Range("J123").Value = 38000
$ws.Range("L123").Value = 38000
$ws.Range("N123").Value = -47800

$ws.Range("H132").Value = 1886653
$ws.Range("I132").Value = 2061.5417
$ws.Range("J132").Value = 24501750
$ws.Range("K132").Value = 6184.625100000001
$ws.Range("L132").Value = 73505250
$ws.Range("M132").Value = -3654.625100000001
$ws.Range("N132").Value = -73510310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 45547440
$ws.Range("I61").Value = 71501510
$ws.Range("J61").Value = 127812.75
$ws.Range("K61").Value = 71501510
$ws.Range("L61").Value = 127812.75
$ws.Range("M61").Value = -71501298
$ws.Range("N61").Value = -128236.75

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H110").Value = 1361.5883
$ws.Range("I110").Value = 662.44446
$ws.Range("K110").Value = 662.44446
$ws.Range("M110").Value = 1382.55554

$ws.Range("H122").Value = 55557056
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -333338230

$ws.Range("H136").Value = 45547440
$ws.Range("I136").Value = 71501510
$ws.Range("J136").Value = 127812.75
$ws.Range("K136").Value = 214504530
$ws.Range("L136").Value = 383438.25
$ws.Range("M136").Value = -214501980
$ws.Range("N136").Value = -388538.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 768.65625
$ws.Range("I20").Value = 565.2381
$ws.Range("K20").Value = 565.2381
$ws.Range("M20").Value = -318.2381

$ws.Range("H107").Value = 1131.3462
$ws.Range("I107").Value = 1086.15
$ws.Range("K107").Value = 1086.15
$ws.Range("M107").Value = 833.8499999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1291.1
$ws.Range("I16").Value = 922.2
$ws.Range("J16").Value = 1660
$ws.Range("K16").Value = 922.2
$ws.Range("L16").Value = 1660
$ws.Range("M16").Value = -635.2
$ws.Range("N16").Value = -2234

$ws.Range("H93").Value = 17803.1
$ws.Range("I93").Value = 17803.1
$ws.Range("K93").Value = 17803.1
$ws.Range("M93").Value = -15931.1

$ws.Range("H105").Value = 858.2222
$ws.Range("I105").Value = 840.5
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 840.5
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 906.5
$ws.Range("N105").Value = -4494

$ws.Range("H107").Value = 467.92856
$ws.Range("I107").Value = 417.8889
$ws.Range("J107").Value = 558
$ws.Range("K107").Value = 417.8889
$ws.Range("L107").Value = 558
$ws.Range("M107").Value = 1502.1111
$ws.Range("N107").Value = -4398

$ws.Range("H113").Value = 1291.1
$ws.Range("I113").Value = 922.2
$ws.Range("J113").Value = 1660
$ws.Range("K113").Value = 922.2
$ws.Range("L113").Value = 1660
$ws.Range("M113").Value = 1247.8
$ws.Range("N113").Value = -6000

$ws.Range("H116").Value = 20742
$ws.Range("J116").Value = 20742
$ws.Range("L116").Value = 20742
$ws.Range("N116").Value = -29920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 2158.3333
$ws.Range("I126").Value = 1333.3334
$ws.Range("J126").Value = 2983.3333
$ws.Range("K126").Value = 4000.0002
$ws.Range("L126").Value = 8949.999899999999
$ws.Range("M126").Value = 939.9998000000001
$ws.Range("N126").Value = -18829.9999

$ws.Range("H131").Value = 950.8169
$ws.Range("I131").Value = 557.6667
$ws.Range("J131").Value = 1007.8871
$ws.Range("K131").Value = 1673.0001
$ws.Range("L131").Value = 3023.6613
$ws.Range("M131").Value = 3366.9999
$ws.Range("N131").Value = -13103.6613

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 622
$ws.Range("N102").ClearContents()

$ws.Range("H113").Value = 1665
$ws.Range("I113").Value = 1400
$ws.Range("K113").Value = 1400
$ws.Range("M113").Value = 770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 825
$ws.Range("I22").Value = 677.2727
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 677.2727
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -382.2727
$ws.Range("N22").Value = -1540

$ws.Range("H27").Value = 825
$ws.Range("I27").Value = 677.2727
$ws.Range("J27").Value = 950
$ws.Range("K27").Value = 677.2727
$ws.Range("L27").Value = 950
$ws.Range("M27").Value = -570.2727
$ws.Range("N27").Value = -1164

$ws.Range("H40").Value = 1900
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1764
$ws.Range("N40").ClearContents()

$ws.Range("H61").Value = 2481.8096
$ws.Range("I61").Value = 2282.8235
$ws.Range("J61").Value = 3327.5
$ws.Range("K61").Value = 2282.8235
$ws.Range("L61").Value = 3327.5
$ws.Range("M61").Value = -2080.8235
$ws.Range("N61").Value = -3731.5

$ws.Range("H68").Value = 1875.6
$ws.Range("I68").Value = 1792.6666
$ws.Range("K68").Value = 1792.6666
$ws.Range("M68").Value = -1043.6666

$ws.Range("H71").Value = 1875.6
$ws.Range("I71").Value = 1792.6666
$ws.Range("K71").Value = 8963.333000000001
$ws.Range("M71").Value = -5219.333000000001

$ws.Range("H113").Value = 2481.8096
$ws.Range("I113").Value = 2282.8235
$ws.Range("J113").Value = 3327.5
$ws.Range("K113").Value = 2282.8235
$ws.Range("L113").Value = 3327.5
$ws.Range("M113").Value = -112.8235
$ws.Range("N113").Value = -7667.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 989.2105
$ws.Range("J113").Value = 1943.5714
$ws.Range("L113").Value = 5830.7142
$ws.Range("N113").Value = -10170.7142

$ws.Range("H117").Value = 20409
$ws.Range("J117").Value = 20409
$ws.Range("L117").Value = 20409
$ws.Range("N117").Value = -29587

$ws.Range("H126").Value = 1189.7
$ws.Range("I126").Value = 1180.875
$ws.Range("J126").Value = 1225
$ws.Range("K126").Value = 3542.625
$ws.Range("L126").Value = 3675
$ws.Range("M126").Value = -1072.625
$ws.Range("N126").Value = -8615

$ws.Range("H132").Value = 70396.10000000001
$ws.Range("I132").Value = 56555.5
$ws.Range("J132").Value = 93044.37
$ws.Range("K132").Value = 169666.5
$ws.Range("L132").Value = 279133.11
$ws.Range("M132").Value = -167136.5
$ws.Range("N132").Value = -284193.11

$ws.Range("H133").Value = 38698.2
$ws.Range("J133").Value = 38698.2
$ws.Range("L133").Value = 38698.2
$ws.Range("N133").Value = -48818.2

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
